$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5467
$ws1.Range("F6").Value = 74
$ws1.Range("F8").Value = 899
$ws1.Range("F9").Value = 140
$ws1.Range("F10").Value = 2429
$ws1.Range("F12").Value = 58
$ws1.Range("F13").Value = 65
$ws1.Range("F14").Value = 2278
$ws1.Range("F15").Value = 162

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 97

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5467
$ws4.Range("F6").Value = 97
$ws4.Range("F7").Value = 74
$ws4.Range("F10").Value = 899
$ws4.Range("F11").Value = 140
$ws4.Range("F12").Value = 2429
$ws4.Range("F14").Value = 58
$ws4.Range("F16").Value = 65
$ws4.Range("F17").Value = 2278
$ws4.Range("F18").Value = 162
